$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.328.72'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '2.520.07'
$ws.Range('E3').Value = '  -5.86%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '574.66'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -4.16%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '169.21'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.57%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -3.56%  '
$ws.Range('D9').Value = '2.519.23'
$ws.Range('E9').Value = '  -5.83%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.163'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').Value = '  -4.21%  '
$ws.Range('E13').Value = '  -4.28%  '
$ws.Range('D14').Value = '2.982.98'
$ws.Range('E14').Value = '  -5.96%  '
$ws.Range('D15').Value = '70.239.87'
$ws.Range('E15').Value = '  -2.22%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000180'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.92%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '24.84'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -5.44%  '
$ws.Range('D18').Value = '2.526.55'
$ws.Range('E18').Value = '  -5.70%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.53'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -5.98%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.54'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -8.24%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '355.99'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -4.45%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '3.92'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -6.08%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.95'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('E24').Value = '  -0.11%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '69.09'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('E26').Value = '  -7.13%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.19'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -6.19%  '
$ws.Range('D28').Value = '2.648.53'
$ws.Range('E28').Value = '  -6.14%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').Value = '0.0₃0908'
$ws.Range('E30').Value = '  -6.82%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.82'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -3.28%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '477.85'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -4.89%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.26'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').Value = '  +0.16%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '156.38'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.66%  '
$ws.Range('E37').Value = '  +4.28%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '18.82'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.36%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.52'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -5.45%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -5.71%  '
$ws.Range('E42').Value = '  -7.85%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.317'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -4.84%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.69'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -6.29%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.40'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -6.51%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '38.26'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -3.09%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '142.91'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -8.64%  '
$ws.Range('E48').Value = '  -5.73%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.522'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -7.13%  '
$ws.Range('E50').Value = '  -7.63%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.594'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.95%  '
